# Insert a new data row at row 347 (Macroferia Regional de Talca - Pepino
# ensalada, weekly price log). This pushes the existing rows 347-424 down to
# 348-425 and adds one brand-new record at row 347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(347).Insert()

$ws.Range("A347").Value = 5
$ws.Range("B347").Value = "Macroferia Regional de Talca"
$ws.Range("C347").Value = "Maule"
$ws.Range("D347").Value = 44711
$ws.Range("E347").Value = 7
$ws.Range("F347").Value = 100112043
$ws.Range("G347").Value = "Pepino ensalada"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 300
$ws.Range("K347").Value = 20000
$ws.Range("L347").Value = 20000
$ws.Range("M347").Value = 20000
$ws.Range("N347").Value = "$/caja 60 unidades"
$ws.Range("O347").Value = "Región de Arica y Parinacota"
$ws.Range("P347").Value = 333
$ws.Range("Q347").Value = 60
$ws.Range("R347").Value = "Hortaliza"
